$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the new "Elvis Presley" test data into row 3 (the row that currently
#    holds the long-form "Elvis Presley (Largest in English src: ...)" name).
$ws.Range("J3").Value = 112.99
$ws.Range("K3").Value = 112.11
$ws.Range("L3").Value = 112.5
$ws.Range("M3").Value = 110.5
$ws.Range("N3").Value = 110.15
$ws.Range("O3").Formula = "=AVERAGE(J3, K3, L3, M3, N3)"
$ws.Range("P3").Formula = "=O3 / 60"
$ws.Range("Q3").Formula = "=O3 / G3"
$ws.Range("R3").Value = "Visual Studio Code"

# 2. Shorten the article name in A3 to just "Elvis Presley".
$ws.Range("A3").Value = "Elvis Presley"

# 3. Center-align every populated cell (this also creates the new cellXfs
#    entry with horizontal="center"). Apply per-row over just the already
#    populated cell span so no new blank cells get materialized.
$ws.Range("A1:R1").HorizontalAlignment = -4108
$ws.Range("A2:R2").HorizontalAlignment = -4108
$ws.Range("A3:R3").HorizontalAlignment = -4108
$ws.Range("A4:B4").HorizontalAlignment = -4108
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("G4:I4").HorizontalAlignment = -4108
$ws.Range("G5:H5").HorizontalAlignment = -4108
$ws.Range("G6:I6").HorizontalAlignment = -4108
$ws.Range("G7:I7").HorizontalAlignment = -4108
$ws.Range("G8:I8").HorizontalAlignment = -4108
$ws.Range("G9:I9").HorizontalAlignment = -4108
$ws.Range("G10:H10").HorizontalAlignment = -4108

# 4. Column width tweaks: B:C pick up the (explicit) default width, R
#    narrows down from its old very-wide size.
$ws.Columns("B:C").ColumnWidth = 9.998697916666666
$ws.Columns("R").ColumnWidth = 16.498697916666668

# 5. Selection / scroll position: back to the default top-left with A4
#    selected (instead of the scrolled-right R9 selection).
$ws.Range("A4").Select()
